$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Volume 30   Number  4" -> "Volume 30   Number  5" in A8 ---
$a8 = $ws.Range("A8")
$a8full = $a8.Value()
$a8ch = $a8.Characters($a8full.Length, 1)
$a8ch.Text = "5"

# --- Update "Report Covering the Week  1/23/2023  Through  1/29/2023" in C9 ---
# "1/23/2023" -> "1/30/2023" ; "1/29/2023" -> "2/5/2023"
$c9 = $ws.Range("C9")
$c9ch1 = $c9.Characters(27, 9)
$c9ch1.Text = "1/30/2023"
$c9ch2 = $c9.Characters(47, 9)
$c9ch2.Text = "2/5/2023"

# --- Data table updates (rows 14-29) ---
$ws.Range("G30").Copy($ws.Range("C14"))
$ws.Range("C14").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 100
$ws.Range("I14").Value = 4
$ws.Range("K14").Value = 100
$ws.Range("M14").Value = 100
$ws.Range("D15").Value = 2
$ws.Range("C30").Copy($ws.Range("F15"))
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = -60
$ws.Range("L15").Value = -50
$ws.Range("N15").Value = -60
$ws.Range("C16").Value = 15
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = 87.5
$ws.Range("F16").Value = 46
$ws.Range("G16").Value = 42
$ws.Range("H16").Value = 9.523809523809
$ws.Range("I16").Value = 62
$ws.Range("J16").Value = 55
$ws.Range("K16").Value = 12.727272727272
$ws.Range("L16").Value = 72.222222222222
$ws.Range("M16").Value = 67.567567567567
$ws.Range("N16").Value = -74.058577405857
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -6.25
$ws.Range("F17").Value = 75
$ws.Range("G17").Value = 60
$ws.Range("H17").Value = 25
$ws.Range("I17").Value = 104
$ws.Range("J17").Value = 75
$ws.Range("K17").Value = 38.666666666666
$ws.Range("L17").Value = 28.395061728395
$ws.Range("M17").Value = 116.666666666667
$ws.Range("N17").Value = -2.803738317757
$ws.Range("C18").Value = 10
$ws.Range("E18").Value = 11.111111111111
$ws.Range("F18").Value = 41
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = 57.692307692307
$ws.Range("I18").Value = 52
$ws.Range("J18").Value = 33
$ws.Range("K18").Value = 57.575757575757
$ws.Range("L18").Value = 126.086956521739
$ws.Range("M18").Value = 33.333333333333
$ws.Range("N18").Value = -79.365079365079
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 53
$ws.Range("G19").Value = 79
$ws.Range("H19").Value = -32.911392405063
$ws.Range("I19").Value = 68
$ws.Range("J19").Value = 91
$ws.Range("K19").Value = -25.274725274725
$ws.Range("L19").Value = -8.108108108108
$ws.Range("M19").Value = 83.783783783783
$ws.Range("N19").Value = 11.475409836065
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 80
$ws.Range("F20").Value = 38
$ws.Range("H20").Value = 52
$ws.Range("I20").Value = 49
$ws.Range("J20").Value = 30
$ws.Range("K20").Value = 63.333333333333
$ws.Range("L20").Value = 444.444444444444
$ws.Range("M20").Value = 308.333333333333
$ws.Range("N20").Value = -72.316384180791
$ws.Range("C21").Value = 70
$ws.Range("D21").Value = 55
$ws.Range("E21").Value = 27.272727272727
$ws.Range("F21").Value = 255
$ws.Range("G21").Value = 238
$ws.Range("H21").Value = 7.142857142857
$ws.Range("I21").Value = 341
$ws.Range("J21").Value = 291
$ws.Range("K21").Value = 17.182130584192
$ws.Range("L21").Value = 50.220264317180
$ws.Range("M21").Value = 90.502793296089
$ws.Range("N21").Value = -59.644970414201
$ws.Range("G30").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 3
$ws.Range("C30").Copy($ws.Range("G22"))
$ws.Range("E30").Copy($ws.Range("H22"))
$ws.Range("I22").Value = 3
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 200
$ws.Range("C30").Copy($ws.Range("D23"))
$ws.Range("E30").Copy($ws.Range("E23"))
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = -25
$ws.Range("I23").Value = 7
$ws.Range("K23").Value = 75
$ws.Range("L23").Value = -30
$ws.Range("M23").Value = 40
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 14.814814814814
$ws.Range("F24").Value = 107
$ws.Range("G24").Value = 99
$ws.Range("H24").Value = 8.080808080808
$ws.Range("I24").Value = 143
$ws.Range("J24").Value = 122
$ws.Range("K24").Value = 17.213114754098
$ws.Range("L24").Value = 4.379562043795
$ws.Range("M24").Value = -2.721088435374
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = -30.434782608695
$ws.Range("F25").Value = 74
$ws.Range("G25").Value = 85
$ws.Range("H25").Value = -12.941176470588
$ws.Range("I25").Value = 92
$ws.Range("J25").Value = 105
$ws.Range("K25").Value = -12.380952380952
$ws.Range("L25").Value = -2.127659574468
$ws.Range("M25").Value = 0
$ws.Range("D26").Value = 5
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 9
$ws.Range("H26").Value = -77.777777777777
$ws.Range("J26").Value = 9
$ws.Range("K26").Value = -55.555555555555
$ws.Range("L26").Value = -33.333333333333
$ws.Range("G30").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 28.571428571428
$ws.Range("I27").Value = 12
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = 71.428571428571
$ws.Range("L27").Value = 140
$ws.Range("C30").Copy($ws.Range("C28"))
$ws.Range("G30").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
$ws.Range("H30").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 6
$ws.Range("H28").Value = 50
$ws.Range("J28").Value = 5
$ws.Range("K28").Value = 40
$ws.Range("L28").Value = 75
$ws.Range("M28").Value = 75
$ws.Range("N28").Value = -30
$ws.Range("C30").Copy($ws.Range("C29"))
$ws.Range("G30").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 1
$ws.Range("H30").Copy($ws.Range("E29"))
$ws.Range("E29").Value = -100
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("J29").Value = 5
$ws.Range("K29").Value = -40
$ws.Range("L29").Value = -25
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -70
